$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 259, shifting existing rows 259:287 down to 261:289
$ws.Range("A259:T260").EntireRow.Insert()

# New row 259 values
$ws.Range("A259").Value = 3
$ws.Range("B259").Value = "Femacal de La Calera"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 44918
$ws.Range("D259").NumberFormat = $ws.Range("D258").NumberFormat
$ws.Range("E259").Value = 5
$ws.Range("F259").Value = "Fruta"
$ws.Range("G259").Value = 100101
$ws.Range("H259").Value = "Berries"
$ws.Range("I259").Value = 100101001
$ws.Range("J259").Value = "Arándano (blue)"
$ws.Range("K259").Value = "Sin especificar"
$ws.Range("L259").Value = "Primera"
$ws.Range("M259").Value = 45
$ws.Range("N259").Value = 4500
$ws.Range("O259").Value = 4500
$ws.Range("P259").Value = 4500
$ws.Range("Q259").Value = "`$/bandeja 2 kilos"
$ws.Range("R259").Value = "Provincia de Curicó"
$ws.Range("S259").Value = 2250
$ws.Range("T259").Value = 2

# New row 260 values
$ws.Range("A260").Value = 3
$ws.Range("B260").Value = "Femacal de La Calera"
$ws.Range("C260").Value = "Coquimbo"
$ws.Range("D260").Value = 44918
$ws.Range("D260").NumberFormat = $ws.Range("D258").NumberFormat
$ws.Range("E260").Value = 5
$ws.Range("F260").Value = "Fruta"
$ws.Range("G260").Value = 100101
$ws.Range("H260").Value = "Berries"
$ws.Range("I260").Value = 100101001
$ws.Range("J260").Value = "Arándano (blue)"
$ws.Range("K260").Value = "Sin especificar"
$ws.Range("L260").Value = "Segunda"
$ws.Range("M260").Value = 48
$ws.Range("N260").Value = 3000
$ws.Range("O260").Value = 3000
$ws.Range("P260").Value = 3000
$ws.Range("Q260").Value = "`$/bandeja 2 kilos"
$ws.Range("R260").Value = "Provincia de Curicó"
$ws.Range("S260").Value = 1500
$ws.Range("T260").Value = 2
